$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "'71.169.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.48%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.870.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.67%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'696.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.86%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'174.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.63%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.868.28"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.66%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.01%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.02%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.19%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'7.19"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -6.50%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -0.56%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000261"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +3.09%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'36.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.88%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.522.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.72%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.891.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.17%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'71.213.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.52%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'17.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.25%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'7.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.22%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -0.31%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'11.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.15%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'498.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +3.79%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.723"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.94%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +3.70%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +1.56%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'10.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +3.16%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'12.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.91%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.15"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.60%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'3.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.77%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.01%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +1.22%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -1.68%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'29.74"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.48%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.182"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.18%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'9.25"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.39%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'3.823.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.78%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.24%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +2.37%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +9.28%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'3.43"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.22%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +8.78%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +0.78%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +0.03%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +0.01%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.000312"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -7.24%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'163.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +2.63%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'49.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.49%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'TheGraph"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'0.303"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.13%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'Bittensor"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'416.47"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +3.52%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -4.00%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -2.32%  "
$ws.Range("E51").Style = "Normal"
